$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 2466
$ws1.Cells.Item(3, 6).Value = 762
$ws1.Cells.Item(4, 6).Value = 250
$ws1.Cells.Item(5, 6).Value = 421
$ws1.Cells.Item(10, 6).Value = 951
$ws1.Cells.Item(14, 6).Value = 76
$ws1.Cells.Item(16, 6).Value = 1110
$ws1.Cells.Item(17, 6).Value = 24471
$ws1.Cells.Item(17, 7).Value = '暂时售罄'
$ws1.Cells.Item(18, 6).Value = 2337
$ws1.Cells.Item(19, 6).Value = 152
$ws1.Cells.Item(20, 6).Value = 366
$ws1.Cells.Item(25, 6).Value = 82
$ws1.Cells.Item(26, 6).Value = 242
$ws1.Cells.Item(28, 6).Value = 78
$ws1.Cells.Item(29, 6).Value = 49
$ws1.Cells.Item(30, 6).Value = 364
$ws1.Cells.Item(32, 6).Value = 449

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(7, 6).Value = 274
$ws2.Cells.Item(8, 6).Value = 167
$ws2.Cells.Item(11, 6).Value = 3650
$ws2.Cells.Item(13, 6).Value = 157

# ---- Sheet 3: 本地生活 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(4, 6).Value = 805

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(4, 6).Value = 2466
$ws4.Cells.Item(5, 6).Value = 805
$ws4.Cells.Item(6, 6).Value = 762
$ws4.Cells.Item(7, 6).Value = 250
$ws4.Cells.Item(8, 6).Value = 421
$ws4.Cells.Item(9, 2).NumberFormat = '@'
$ws4.Cells.Item(9, 2).Value = '2024-06-28'
$ws4.Cells.Item(9, 3).Value = '广州·2024日本电音大神Ken Arai中国巡演'
$ws4.Cells.Item(9, 4).Value = '恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）'
$ws4.Cells.Item(9, 5).Value = '2024.06.28 20:00-06.28 21:30'
$ws4.Cells.Item(9, 6).Value = 1
$ws4.Cells.Item(9, 7).Value = 380
$ws4.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87429'
$ws4.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/hTVOnzgo1717571680673.jpeg'
$ws4.Cells.Item(10, 3).Value = '广州·奥斯卡·罗曼耶卓（O叔）钢琴独奏音乐会'
$ws4.Cells.Item(10, 4).Value = '晴波路33号 广州星海音乐厅'
$ws4.Cells.Item(10, 6).Value = 189
$ws4.Cells.Item(10, 7).Value = 480
$ws4.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84545'
$ws4.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/XK8EYxGv1712890578712.jpeg'
$ws4.Cells.Item(11, 2).NumberFormat = '@'
$ws4.Cells.Item(11, 2).Value = '2024-06-29'
$ws4.Cells.Item(11, 3).Value = '广州·《千与千寻》宫崎骏&久石让经典动漫歌曲视听音乐会'
$ws4.Cells.Item(11, 4).Value = '香雪大道西3号 广州科学城会议中心'
$ws4.Cells.Item(11, 5).Value = '2024.06.29 19:30-06.29 21:00'
$ws4.Cells.Item(11, 6).Value = 9
$ws4.Cells.Item(11, 7).Value = 76.5
$ws4.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87328'
$ws4.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/QMBx7g3M1718267752079.png'
$ws4.Cells.Item(12, 3).Value = '广州·《海上钢琴师》经典电影作品大型交响音乐会'
$ws4.Cells.Item(12, 4).Value = '东风中路299号 广州中山纪念堂'
$ws4.Cells.Item(12, 5).Value = '2024.06.29 20:00-06.29 21:40'
$ws4.Cells.Item(12, 6).Value = 93
$ws4.Cells.Item(12, 7).Value = 75
$ws4.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84162'
$ws4.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/OnWieQKe1712742593534.jpeg'
$ws4.Cells.Item(13, 3).Value = '广州·掠空纪演唱会'
$ws4.Cells.Item(13, 4).Value = '和平路1号 良仓新造创意园'
$ws4.Cells.Item(13, 5).Value = '2024.06.29 19:00-06.29 22:00'
$ws4.Cells.Item(13, 6).Value = 273
$ws4.Cells.Item(13, 7).Value = 290
$ws4.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86991'
$ws4.Cells.Item(13, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/joOOx0Cr1717669820698.jpeg'
$ws4.Cells.Item(14, 2).NumberFormat = '@'
$ws4.Cells.Item(14, 2).Value = '2024-07-01'
$ws4.Cells.Item(14, 3).Value = '广州·NIJISANJI EN 官方授权主题店'
$ws4.Cells.Item(14, 4).Value = '天河路299号 时尚天河商业广场'
$ws4.Cells.Item(14, 5).Value = '2024.07.01 00:00-07.15 23:59'
$ws4.Cells.Item(14, 6).Value = 246
$ws4.Cells.Item(14, 7).Value = 30
$ws4.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86338'
$ws4.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202405/dB7yQHnF1716795983429.jpeg'
$ws4.Cells.Item(15, 2).NumberFormat = '@'
$ws4.Cells.Item(15, 2).Value = '2024-07-06'
$ws4.Cells.Item(15, 3).Value = '广州·OVO动漫展'
$ws4.Cells.Item(15, 4).Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws4.Cells.Item(15, 5).Value = '2024.07.06 10:00-07.06 17:00'
$ws4.Cells.Item(15, 6).Value = 913
$ws4.Cells.Item(15, 7).Value = 45
$ws4.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85706'
$ws4.Cells.Item(15, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/NP2wHpFI1715762116387.png'
$ws4.Cells.Item(16, 3).Value = '广州·火影忍者only'
$ws4.Cells.Item(16, 6).Value = 573
$ws4.Cells.Item(16, 7).Value = 60
$ws4.Cells.Item(16, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85704'
$ws4.Cells.Item(16, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/lKOROXve1715763433389.jpeg'
$ws4.Cells.Item(17, 3).Value = '广州·第五人格ONLY2.0'
$ws4.Cells.Item(17, 4).Value = '洛浦街夏滘西环路1号(厦滘地铁站A口步行290米) 厦喾岭南电商园会展中心'
$ws4.Cells.Item(17, 6).Value = 949
$ws4.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86136'
$ws4.Cells.Item(17, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/JBAZGW9P1716363092370.jpeg'
$ws4.Cells.Item(18, 2).NumberFormat = '@'
$ws4.Cells.Item(18, 2).Value = '2024-07-13'
$ws4.Cells.Item(18, 3).Value = '广州·特摄FansMeetup'
$ws4.Cells.Item(18, 4).Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws4.Cells.Item(18, 5).Value = '2024.07.13 10:00-07.13 19:00'
$ws4.Cells.Item(18, 6).Value = 132
$ws4.Cells.Item(18, 7).Value = 69.9
$ws4.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87031'
$ws4.Cells.Item(18, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/9ffC9a8n1717578946827.jpeg'
$ws4.Cells.Item(19, 3).Value = '广州·第5人格only3.0联动特别篇'
$ws4.Cells.Item(19, 4).Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws4.Cells.Item(19, 5).Value = '2024.07.13 10:00-07.13 17:00'
$ws4.Cells.Item(19, 6).Value = 443
$ws4.Cells.Item(19, 7).Value = 60
$ws4.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86740'
$ws4.Cells.Item(19, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/mwlJqj0o1717149700846.jpeg'
$ws4.Cells.Item(20, 2).NumberFormat = '@'
$ws4.Cells.Item(20, 2).Value = '2024-07-14'
$ws4.Cells.Item(20, 3).Value = '广州·OOPS 7th'
$ws4.Cells.Item(20, 4).Value = '流花街道流花路119号越秀公园站B2、C出口 广州越秀国际会议中心'
$ws4.Cells.Item(20, 5).Value = '2024.07.14 09:30-07.15 17:00'
$ws4.Cells.Item(20, 6).Value = 74
$ws4.Cells.Item(20, 7).Value = 20
$ws4.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87550'
$ws4.Cells.Item(20, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/Qi8gB0Bi1715922859908.png'
$ws4.Cells.Item(21, 3).Value = '广州·幻毛纪AnimalFurryOnly'
$ws4.Cells.Item(21, 4).Value = '芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋'
$ws4.Cells.Item(21, 5).Value = '2024.07.14 10:00-07.14 19:00'
$ws4.Cells.Item(21, 6).Value = 35
$ws4.Cells.Item(21, 7).Value = 68.8
$ws4.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87273'
$ws4.Cells.Item(21, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/9z1DMHsl1718181280279.png'
$ws4.Cells.Item(22, 3).Value = '广州·火影only'
$ws4.Cells.Item(22, 4).Value = '人和镇蚌湖清河大街168号 人和园'
$ws4.Cells.Item(22, 5).Value = '2024.07.14 09:30-07.14 17:30'
$ws4.Cells.Item(22, 6).Value = 1106
$ws4.Cells.Item(22, 7).Value = 78
$ws4.Cells.Item(22, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84815'
$ws4.Cells.Item(22, 9).Value = '//i2.hdslb.com/bfs/openplatform/202404/QLOhW4Nr1714384036670.png'
$ws4.Cells.Item(23, 2).NumberFormat = '@'
$ws4.Cells.Item(23, 2).Value = '2024-07-19'
$ws4.Cells.Item(23, 3).Value = '广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园'
$ws4.Cells.Item(23, 4).Value = '新港东路1000号 保利世贸博览馆'
$ws4.Cells.Item(23, 5).Value = '2024.07.19 09:00-07.22 17:00'
$ws4.Cells.Item(23, 6).Value = 24471
$ws4.Cells.Item(23, 7).Value = '暂时售罄'
$ws4.Cells.Item(23, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87210'
$ws4.Cells.Item(23, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg'
$ws4.Cells.Item(24, 6).Value = 24471
$ws4.Cells.Item(24, 7).Value = '暂时售罄'
$ws4.Cells.Item(28, 6).Value = 157
$ws4.Cells.Item(30, 6).Value = 2337
$ws4.Cells.Item(31, 6).Value = 152
$ws4.Cells.Item(33, 6).Value = 366
$ws4.Cells.Item(38, 6).Value = 242
$ws4.Cells.Item(41, 6).Value = 78
$ws4.Cells.Item(42, 6).Value = 49
$ws4.Cells.Item(46, 6).Value = 449

